$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text values (Price column) must be forced to
# Text format first so Excel keeps them as strings (matching the original
# inlineStr storage) instead of converting them to numbers.
$numericTextCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D40", "D41", "D42", "D43", "D45", "D49")
foreach ($c in $numericTextCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "244.05"
$ws.Range("D3").Value = "22.99"
$ws.Range("D4").Value = "5.405"
$ws.Range("D5").Value = "0.05971"
$ws.Range("D6").Value = "3.460"
$ws.Range("D7").Value = "6.524"
$ws.Range("D8").Value = "0.8140"
$ws.Range("D9").Value = "0.9156"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01127"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1408"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07406"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03237"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03088"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09359"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.854"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001558"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04667"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.006087"
$ws.Range("D21").Value = "0.0009827"
$ws.Range("D22").Value = "0.00007801"
$ws.Range("D23").Value = "3.614"
$ws.Range("D40").Value = "0.03926"
$ws.Range("D41").Value = "0.006241"
$ws.Range("D42").Value = "0.1074"
$ws.Range("D43").Value = "0.003000"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "0.00005248"
$ws.Range("D49").Value = "0.002288"

# Restore the default (no explicit number format) style so the cells stay
# visually/structurally identical to the rest of the sheet.
foreach ($c in $numericTextCells) {
    $ws.Range($c).Style = "Normal"
}
